$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ST3MR-5 - As a shop user i want to have the option Other in the Categories module in order to buy products - 1 defect", "ST3MR-8 - Verify that the Other Tools products sort is visible", "eCommerce - Other Tools  - FAIL", "ST3MR-10 - In Progress - The Sort products dropdown box don't show nothing"),
    @("ST3MR-5 - As a shop user i want to have the option Other in the Categories module in order to buy products", "ST3MR-20 - Verify the Categories module on Homepage", "eCommerce - Other Tools  - PASS", "None"),
    @("ST3MR-5 - As a shop user i want to have the option Other in the Categories module in order to buy products", "ST3MR-21 - Verify the Other Tools on Categories module", "eCommerce - Other Tools  - PASS", "None"),
    @("ST3MR-3 - As a shop user i want to have the option Power Tools in the Categories module in order to buy products", "ST3MR-9 - Verify that the Power Tools products filters is visible", "eCommerce - Power Tools  - PASS", "None"),
    @("ST3MR-3 - As a shop user i want to have the option Power Tools in the Categories module in order to buy products", "ST3MR-18 - Verify the Categories module on Homepage", "eCommerce - Power Tools  - PASS", "None"),
    @("ST3MR-3 - As a shop user i want to have the option Power Tools in the Categories module in order to buy products", "ST3MR-19 - Verify the Power Tools on Categories module", "eCommerce - Power Tools  - PASS", "None"),
    @("ST3MR-2 - As a shop user i want to have the option Hand Tools in the Categories module in order to buy products", "ST3MR-13 - Verify that the user can add products to cart", "eCommerce - Hand Tools  - PASS", "None"),
    @("ST3MR-2 - As a shop user i want to have the option Hand Tools in the Categories module in order to buy products - 1 defect", "ST3MR-16 - Verify that the user can add products to favourites", "eCommerce - Hand Tools  - FAIL", "ST3MR-17 - In Progress - The Add to Favourites command doesn't work  - step"),
    @("ST3MR-2 - As a shop user i want to have the option Hand Tools in the Categories module in order to buy products", "ST3MR-7 - Verify all products from Hand Tools category and products detail page", "eCommerce - Hand Tools  - PASS", "None"),
    @("ST3MR-2 - As a shop user i want to have the option Hand Tools in the Categories module in order to buy products", "ST3MR-24 - Verify the Hand Tool product can be selected", "eCommerce - Hand Tools  - PASS", "None"),
    @("ST3MR-2 - As a shop user i want to have the option Hand Tools in the Categories module in order to buy products", "ST3MR-22 - Verify the Categories module on Homepage", "eCommerce - Hand Tools  - PASS", "None"),
    @("ST3MR-2 - As a shop user i want to have the option Hand Tools in the Categories module in order to buy products", "ST3MR-23 - Verify the Hand Tools on Categories module", "eCommerce - Hand Tools  - PASS", "None")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
